# Update countries & provincias Spain
# Applies the data refresh captured in the commit diff:
#  - Swaps the row content (country name + stats) for Panama/Kazajistan/Belgica (rows 36-38)
#  - Swaps the row content (country name + stats) for Islas Malvinas/Montserrat (rows 215-216)
#  - Refreshes the case-count statistics for several countries (rows 22, 50, 107, 132)
#  - Updates the "last updated" timestamp string in cell A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        [string]$Label,
        [double]$B,
        [double]$C,
        [double]$D,
        [double]$E,
        [double]$F,
        [double]$G,
        [double]$H
    )
    $ws.Cells.Item($Row, 1).Value = $Label
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
}

# Row 22 - Pakistan: refreshed totals
Set-Row 22 "Pakistan" 309015 798 294740 7831 0 7 6444

# Rows 36-38: country order/data rotates (Panama -> Belgica -> Kazajistan -> Panama)
Set-Row 36 "Belgica" 108768 1881 19123 79680 0 6 9965
Set-Row 37 "Panama" 108726 0 85494 20935 0 0 2297
Set-Row 38 "Kazajistan" 107590 61 102360 3531 0 0 1699

# Row 50 - Honduras: refreshed totals
Set-Row 50 "Honduras" 73193 518 24580 46364 0 27 2249

# Row 107 - Birmania: refreshed totals
Set-Row 107 "Birmania" 8515 171 2381 5979 0 5 155

# Row 132 - Trinidad yTobago: refreshed totals
Set-Row 132 "Trinidad yTobago" 4235 0 2047 2121 0 0 67

# Rows 215-216: Islas Malvinas / Montserrat swap places
Set-Row 215 "Montserrat" 13 0 12 0 0 0 1
Set-Row 216 "Islas Malvinas" 13 0 13 0 0 0 0

# Update the "last refreshed" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Septiembre de 2020 a las 05:12"
